$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = 1294912.0
$ws.Range("D6").Value = 2450768.0
